$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Theodore Open Display Cabinet" row (row 1) is removed entirely, which
# shifts "Rhodes Folding Book Shelf" / "Austen Bookshelf/Display Unit" up one
# row. A new row 3 is then added for "Rhodes Folding Book Shelf" at its
# updated (modified) price.
$ws.Rows("1:1").Delete()

$ws.Range("A3").Value = "Rhodes Folding Book Shelf"
$ws.Range("B3").Value = "₹12,039"

$ws.Columns("A:A").AutoFit()
